# ValueSet-fr-editorial-status.xlsx
# - bump the "Date" metadata value to the new commit timestamp
# - set the (previously empty) "Jurisdiction" value to "FRANCE"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 8: Date | 2025-07-11T12:24:54+00:00  ->  2025-07-11T12:29:53+00:00
$ws.Cells.Item(8, 2).Value = "2025-07-11T12:29:53+00:00"

# Row 11: Jurisdiction | (empty)  ->  FRANCE
$ws.Cells.Item(11, 2).Value = "FRANCE"
